$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): argument / predicate labels
$ws.Range("B1").Value = "arg1"
$ws.Range("C1").Value = "arg2"
$ws.Range("D1").Value = "arg3"
$ws.Range("E1").Value = "arg4"
$ws.Range("F1").Value = "pred1"
$ws.Range("G1").Value = "pred2"
$ws.Range("H1").Value = "pred3"
$ws.Range("I1").Value = "pred4"

# Data row (row 2)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "to the name"
$ws.Range("C2").Value = "their name"
$ws.Range("D2").Value = "common"
$ws.Range("E2").Value = "it"
$ws.Range("F2").Value = "On the other hand , things are said to be named Univocally which have both the name and the definition answering to the name in common ."
$ws.Range("G2").Value = "Things are said to be named Derivatively , which derive their name from some other name , but differ from it in termination ."
$ws.Range("H2").Value = "both the name answering to the name in common"
$ws.Range("I2").Value = "to be which derive their name from some other name"

# Build the bold / bordered / centered-top-aligned format on a single cell
# first, then fan it out with a format-only paste so the workbook ends up
# with exactly one extra style entry (matches how Excel itself collapses
# repeated identical formatting into a single shared cell style record).
$model = $ws.Range("B1")
$model.Font.Bold = $true
$model.Borders.LineStyle = 1
$model.HorizontalAlignment = -4108
$model.VerticalAlignment = -4160

$model.Copy()
$ws.Range("C1:I1").PasteSpecial(-4122)
$model.Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
